# Flow of the project established
# Update the "Before Removing outliers" continuous-variable summary stats
# for the first data row (HRLYEARN) and let the bestFit-style column
# widths follow the new, narrower values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated summary statistics (row 2: D=Mean, E=Median, F=Min, G=Max, H=StdDev) ---
$ws.Range("D2").Value = 3.51
$ws.Range("E2").Value = 3.49
$ws.Range("F2").Value = 1.75
$ws.Range("G2").Value = 5.33
$ws.Range("H2").Value = 0.44

# --- Column widths shrink to fit the new, shorter numbers ---
# (columns 4 and 6 keep their existing best-fit width unchanged)
$ws.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws.Columns.Item(5).ColumnWidth = 6.333333333333333
$ws.Columns.Item(7).ColumnWidth = 4.333333333333333
$ws.Columns.Item(8).ColumnWidth = 16.333333333333332
